# Fruta / hortaliza, semanal
# Insert a new weekly record before the existing row 349 (Femacal de La
# Calera - Espinaca), shifting rows 349-366 down to 350-367, and fill the
# new row 349 with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 349; everything currently at/after 349
# shifts down by one (349->350, ..., 366->367).
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with the new weekly entry.
$ws.Range("A349").Value = 3
$ws.Range("B349").Value = "Femacal de La Calera"
$ws.Range("C349").Value = "Coquimbo"
$ws.Range("D349").Value = 44753
$ws.Range("E349").Value = 5
$ws.Range("F349").Value = 100112012
$ws.Range("G349").Value = "Espinaca"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 230
$ws.Range("K349").Value = 4000
$ws.Range("L349").Value = 4500
$ws.Range("M349").Value = 4261
$ws.Range("N349").Value = '$/docena de atados (3 kilos)'
$ws.Range("O349").Value = "Provincia de Quillota"
$ws.Range("P349").Value = 1420
$ws.Range("Q349").Value = 3
$ws.Range("R349").Value = "Hortaliza"
